$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: (Ryanair flight FR2678, Sunday Jan 15)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Sunday, Jan 15"
$ws.Range("C11").Value = "4:10 PM"
$ws.Range("D11").Value = "FR2678"
$ws.Range("E11").Value = "London"
$ws.Range("F11").Value = "(STN)"
$ws.Range("G11").Value = "Ryanair "
$ws.Range("H11").Value = "B738"
$ws.Range("I11").Value = "(EI-DYN)"
$ws.Range("J11").Value = "4:24 PM"
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = "0 hours, 14 minutes"
$ws.Range("M11").Value = ""

# New row 12: (LOT flight LO3507, Sunday Jan 15)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Sunday, Jan 15"
$ws.Range("C12").Value = "5:55 PM"
$ws.Range("D12").Value = "LO3507"
$ws.Range("E12").Value = "Krakow"
$ws.Range("F12").Value = "(KRK)"
$ws.Range("G12").Value = "LOT "
$ws.Range("H12").Value = "E75S"
$ws.Range("I12").Value = "(SP-LIA)"
$ws.Range("J12").Value = "5:15 PM"
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = "0 hours, -40 minutes"
$ws.Range("M12").Value = ""
